# Auto-generated edit script: update cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.837.16"
$ws.Range("E2").Value = "  +0.13%  "
# Row 3
$ws.Range("D3").Value = "2.665.38"
$ws.Range("E3").Value = "  -0.53%  "
# Row 4
$ws.Range("E4").Value = "  +0.03%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.65"
$ws.Range("E5").Value = "  -0.62%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.19"
$ws.Range("E6").Value = "  +0.82%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.653"
$ws.Range("E7").Value = "  +4.49%  "
# Row 9
$ws.Range("E9").Value = "  -2.79%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.404"
$ws.Range("E10").Value = "  +0.51%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.87"
$ws.Range("E11").Value = "  -0.11%  "
# Row 12
$ws.Range("E12").Value = "  +1.53%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.05"
$ws.Range("E13").Value = "  -1.24%  "
# Row 14
$ws.Range("E14").Value = "  -2.08%  "
# Row 15
$ws.Range("D15").Value = "3.143.04"
$ws.Range("E15").Value = "  -0.61%  "
# Row 16
$ws.Range("D16").Value = "65.702.31"
$ws.Range("E16").Value = "  +0.12%  "
# Row 17
$ws.Range("D17").Value = "2.629.24"
$ws.Range("E17").Value = "  -1.87%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.64"
$ws.Range("E18").Value = "  -2.07%  "
# Row 19
$ws.Range("E19").Value = "  +0.05%  "
# Row 20
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.49"
$ws.Range("E20").Value = "  -1.28%  "
# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.43"
$ws.Range("E21").Value = "  -0.29%  "
# Row 22
$ws.Range("E22").Value = "  -0.08%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.82"
$ws.Range("E23").Value = "  +0.11%  "
# Row 24
$ws.Range("E24").Value = "  +12.12%  "
# Row 25
$ws.Range("E25").Value = "  +0.71%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.74"
$ws.Range("E26").Value = "  +0.79%  "
# Row 27
$ws.Range("E27").Value = "  +1.48%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "568.56"
$ws.Range("E28").Value = "  +7.02%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.23"
$ws.Range("E29").Value = "  +1.69%  "
# Row 30
$ws.Range("E30").Value = "  -2.81%  "
# Row 31
$ws.Range("E31").Value = "  -0.02%  "
# Row 32
$ws.Range("E32").Value = "  -0.39%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.84"
$ws.Range("E33").Value = "  +4.39%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.71"
$ws.Range("E34").Value = "  +3.41%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.62"
$ws.Range("E35").Value = "  +2.10%  "
# Row 36
$ws.Range("E36").Value = "  -0.64%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.65"
$ws.Range("E37").Value = "  +0.40%  "
# Row 38
$ws.Range("E38").Value = "  +0.01%  "
# Row 39
$ws.Range("E39").Value = "  +0.51%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.72"
$ws.Range("E40").Value = "  -2.24%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "161.54"
$ws.Range("E41").Value = "  -1.95%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.11"
$ws.Range("E42").Value = "  -0.98%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0623"
$ws.Range("E43").Value = "  +1.99%  "
# Row 44
$ws.Range("E44").Value = "  -0.63%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.10"
$ws.Range("E45").Value = "  +0.72%  "
# Row 46
$ws.Range("E46").Value = "  +0.28%  "
# Row 47
$ws.Range("E47").Value = "  +0.01%  "
# Row 48
$ws.Range("E48").Value = "  +1.47%  "
# Row 49
$ws.Range("E49").Value = "  -1.79%  "
# Row 50
$ws.Range("D50").Value = "0.0₆0246"
$ws.Range("E50").Value = "  -4.70%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.818"
$ws.Range("E51").Value = "  +0.02%  "
